# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.404.07'
$ws.Range('E2').Value = '  +4.88%  '
$ws.Range('D3').Value = '2.750.27'
$ws.Range('E3').Value = '  +4.78%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'115.67"
$ws.Range('E5').Value = '  +3.96%  '
$ws.Range('D6').Value = "'332.36"
$ws.Range('E6').Value = '  +3.08%  '
$ws.Range('D7').Value = "'0.539"
$ws.Range('E7').Value = '  +2.65%  '
$ws.Range('D8').Value = "'1.00"
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +5.99%  '
$ws.Range('D10').Value = "'41.52"
$ws.Range('E10').Value = '  +4.79%  '
$ws.Range('E11').Value = '  +6.07%  '
$ws.Range('D12').Value = "'20.21"
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('D14').Value = "'7.61"
$ws.Range('E14').Value = '  +5.12%  '
$ws.Range('D15').Value = '3.180.38'
$ws.Range('E15').Value = '  +4.71%  '
$ws.Range('D16').Value = '2.750.87'
$ws.Range('E16').Value = '  +4.53%  '
$ws.Range('D17').Value = "'0.883"
$ws.Range('E17').Value = '  +3.24%  '
$ws.Range('D18').Value = '51.452.25'
$ws.Range('E18').Value = '  +4.98%  '
$ws.Range('D19').Value = "'3.23"
$ws.Range('E19').Value = '  +7.39%  '
$ws.Range('D20').Value = "'13.45"
$ws.Range('E20').Value = '  +4.22%  '
$ws.Range('E21').Value = '  +2.57%  '
$ws.Range('E22').Value = '  +3.37%  '
$ws.Range('D23').Value = "'277.73"
$ws.Range('E23').Value = '  +3.30%  '
$ws.Range('D24').Value = "'69.52"
$ws.Range('E25').Value = '  +4.26%  '
$ws.Range('D26').Value = "'26.79"
$ws.Range('E26').Value = '  +2.77%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('E29').Value = '  -0.67%  '
$ws.Range('E30').Value = '  +2.18%  '
$ws.Range('D31').Value = "'35.09"
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('D33').Value = "'5.54"
$ws.Range('E33').Value = '  +1.37%  '
$ws.Range('D34').Value = "'0.0825"
$ws.Range('E34').Value = '  +3.38%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').Value = "'19.07"
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = "'2.09"
$ws.Range('E37').Value = '  +2.75%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'4.99"
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').Value = "'3.27"
$ws.Range('E39').Value = '  +5.19%  '
$ws.Range('D40').Value = "'126.93"
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('D41').Value = "'23.15"
$ws.Range('E41').Value = '  +4.45%  '
$ws.Range('E42').Value = '  +8.20%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = "'0.0345"
$ws.Range('E43').Value = '  +9.07%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = "'0.114"
$ws.Range('E44').Value = '  +2.97%  '
$ws.Range('D45').Value = "'2.44"
$ws.Range('E45').Value = '  +13.18%  '
$ws.Range('D46').Value = '2.091.49'
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('D47').Value = "'3.32"
$ws.Range('E47').Value = '  +3.82%  '
$ws.Range('E48').Value = '  +4.57%  '
$ws.Range('E49').Value = '  +6.37%  '
$ws.Range('D50').Value = "'8.98"
$ws.Range('E50').Value = '  +1.18%  '
$ws.Range('D51').Value = "'59.88"
$ws.Range('E51').Value = '  +2.23%  '
